{"js": "// COB-27: user manual update for configure number of recent files\n//\n// Insert a new bulleted/numbered list item right after the paragraph\n// \"A menu where the user can select a visual theme for the cobbler\n// application, including a dark mode.\" in the Application Settings\n// section, describing the new \"number of recent files\" setting.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"A menu where the user can select a visual theme for the cobbler application, including a dark mode.\";\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\n    \"Could not locate the anchor paragraph (visual theme setting) to insert the new list item after.\"\n  );\n}\n\n// insertParagraph(\"After\") clones the source paragraph's formatting\n// (pStyle \"ListParagraph\" + numPr ilvl/numId), so the new bullet lands\n// in the same list as its neighbours.\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"A menu that allows user to specify the number of recent files to track and display in the Recent Files menu. \",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# COB-27: user manual update for configure number of recent files\n#\n# Insert a new bulleted/numbered list item right after the paragraph\n# \"A menu where the user can select a visual theme for the cobbler\n# application, including a dark mode.\" in the Application Settings\n# section, describing the new \"number of recent files\" setting.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"A menu where the user can select a visual theme for the cobbler application, including a dark mode.\"\n\n$anchorParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $anchorText) {\n        $anchorParagraph = $p\n        break\n    }\n}\n\nif ($anchorParagraph -eq $null) {\n    throw \"Could not locate the anchor paragraph (visual theme setting) to insert the new list item after.\"\n}\n\n# InsertParagraphAfter() clones the source paragraph's formatting\n# (pStyle \"ListParagraph\" + numPr ilvl/numId), so the new bullet lands\n# in the same list as its neighbours. Like real Word COM, it returns\n# nothing, so re-fetch the freshly created paragraph via .Next() and\n# set its text.\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph = $anchorParagraph.Next()\n$newParagraph.Range.Text = \"A menu that allows user to specify the number of recent files to track and display in the Recent Files menu. \"\n"}
